$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the single registered athlete's name / e-mail (row 2)
$ws.Range("B2").Value = "athlete290"
$ws.Range("C2").Value = "athlete290@smashhub.com"

# Add a new (still empty) member row, pre-formatted with the Hyperlink
# style so the Email column is ready for the next entry.
$ws.Range("C3").Style = "Hyperlink"

# Move the selection roughly to where the author left it.
$ws.Range("H9:H10").Select()
